$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (exhibitions)
$wsExpo = $wb.Worksheets.Item(1)
$wsExpo.Range("F4").Value  = 867
$wsExpo.Range("F5").Value  = 31
$wsExpo.Range("F6").Value  = 326
$wsExpo.Range("F7").Value  = 10303
$wsExpo.Range("F8").Value  = 83
$wsExpo.Range("F18").Value = 810
$wsExpo.Range("F20").Value = 99

# Sheet 2: 演出 (performances)
$wsShow = $wb.Worksheets.Item(2)
$wsShow.Range("F2").Value = 10

# Sheet 4: 全部类型 (all types combined)
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F4").Value  = 867
$wsAll.Range("F5").Value  = 31
$wsAll.Range("F6").Value  = 326
$wsAll.Range("F7").Value  = 10303
$wsAll.Range("F8").Value  = 83
$wsAll.Range("F9").Value  = 0
$wsAll.Range("F18").Value = 810
$wsAll.Range("F20").Value = 99
$wsAll.Range("F21").Value = 10
